# Auto-generated Excel COM-interop script applying the Garuda_Profits price update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value2 = 2639.3
$ws.Range("I29").Value2 = 863.3333
$ws.Range("J29").Value2 = 3400.4285
$ws.Range("K29").Value2 = 2589.9999
$ws.Range("L29").Value2 = 10201.2855
$ws.Range("M29").Value2 = -2308.9999
$ws.Range("N29").Value2 = -10763.2855
$ws.Range("H107").Value2 = 332.83334
$ws.Range("I107").Value2 = 345.53845
$ws.Range("J107").Value2 = 299.8
$ws.Range("K107").Value2 = 345.53845
$ws.Range("L107").Value2 = 299.8
$ws.Range("M107").Value2 = 1574.46155
$ws.Range("N107").Value2 = -4139.8
$ws.Range("H116").Value2 = 1833.3334
$ws.Range("I116").Value2 = 1500
$ws.Range("K116").Value2 = 1500
$ws.Range("M116").Value2 = 1942
$ws.Range("H132").Value2 = 3402649.8
$ws.Range("I132").Value2 = 3572682.2
$ws.Range("J132").Value2 = 2000
$ws.Range("K132").Value2 = 10718046.6
$ws.Range("L132").Value2 = 6000
$ws.Range("M132").Value2 = -10715516.6
$ws.Range("N132").Value2 = -11060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 76924360
$ws.Range("I2").Value2 = 100000820
$ws.Range("K2").Value2 = 100000820
$ws.Range("M2").Value2 = -100000707
$ws.Range("H45").Value2 = 19609300
$ws.Range("I45").Value2 = 55556724
$ws.Range("J45").Value2 = 1613.6364
$ws.Range("K45").Value2 = 55556724
$ws.Range("L45").Value2 = 1613.6364
$ws.Range("M45").Value2 = -55556347
$ws.Range("N45").Value2 = -2367.6364
$ws.Range("H61").Value2 = 2090.0625
$ws.Range("I61").Value2 = 1226.8182
$ws.Range("K61").Value2 = 1226.8182
$ws.Range("M61").Value2 = -1014.8182
$ws.Range("H63").Value2 = 2001599.8
$ws.Range("I63").Value2 = 2001599.8
$ws.Range("J63").Value2 = 0
$ws.Range("K63").Value2 = 2001599.8
$ws.Range("L63").Value2 = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value2 = -2000913.8
$ws.Range("H66").Value2 = 2001599.8
$ws.Range("I66").Value2 = 2001599.8
$ws.Range("J66").Value2 = 0
$ws.Range("K66").Value2 = 10007999
$ws.Range("L66").Value2 = 0
$ws.Range("M66").Value2 = -10004567
$ws.Range("N66").ClearContents()
$ws.Range("H92").Value2 = 49553.8
$ws.Range("J92").Value2 = 49553.8
$ws.Range("L92").Value2 = 49553.8
$ws.Range("N92").Value2 = -54545.8
$ws.Range("H110").Value2 = 2189.111
$ws.Range("J110").Value2 = 0
$ws.Range("L110").Value2 = 0
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value2 = 76924360
$ws.Range("I116").Value2 = 100000820
$ws.Range("K116").Value2 = 100000820
$ws.Range("M116").Value2 = -99998526
$ws.Range("H136").Value2 = 2090.0625
$ws.Range("I136").Value2 = 1226.8182
$ws.Range("K136").Value2 = 3680.4546
$ws.Range("M136").Value2 = -1130.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 76924360
$ws.Range("I3").Value2 = 100000820
$ws.Range("K3").Value2 = 100000820
$ws.Range("M3").Value2 = -100000706
$ws.Range("H134").Value2 = 69028.266
$ws.Range("I134").Value2 = 144403.42
$ws.Range("K134").Value2 = 433210.26
$ws.Range("M134").Value2 = -430675.26

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2927071.2
$ws.Range("I31").Value2 = 2273.5356
$ws.Range("J31").Value2 = 5751014
$ws.Range("K31").Value2 = 2273.5356
$ws.Range("L31").Value2 = 5751014
$ws.Range("M31").Value2 = -1978.5356
$ws.Range("N31").Value2 = -5751604
$ws.Range("H34").Value2 = 2927071.2
$ws.Range("I34").Value2 = 2273.5356
$ws.Range("J34").Value2 = 5751014
$ws.Range("K34").Value2 = 2273.5356
$ws.Range("L34").Value2 = 5751014
$ws.Range("M34").Value2 = -2071.5356
$ws.Range("N34").Value2 = -5751418
$ws.Range("H58").Value2 = 2050
$ws.Range("I58").Value2 = 2300
$ws.Range("J58").Value2 = 1425
$ws.Range("K58").Value2 = 2300
$ws.Range("L58").Value2 = 1425
$ws.Range("M58").Value2 = -2097
$ws.Range("N58").Value2 = -1831
$ws.Range("H99").Value2 = 2520.85
$ws.Range("I99").Value2 = 2525.5715
$ws.Range("J99").Value2 = 2518.3076
$ws.Range("K99").Value2 = 2525.5715
$ws.Range("L99").Value2 = 2518.3076
$ws.Range("M99").Value2 = -1027.5715
$ws.Range("N99").Value2 = -5514.3076
$ws.Range("H105").Value2 = 701.63635
$ws.Range("I105").Value2 = 575.4545000000001
$ws.Range("J105").Value2 = 827.8182
$ws.Range("K105").Value2 = 575.4545000000001
$ws.Range("L105").Value2 = 827.8182
$ws.Range("M105").Value2 = 1171.5455
$ws.Range("N105").Value2 = -4321.8182
$ws.Range("H108").Value2 = 13666.667
$ws.Range("J108").Value2 = 13666.667
$ws.Range("L108").Value2 = 13666.667
$ws.Range("N108").Value2 = -21346.667
$ws.Range("H126").Value2 = 2520.85
$ws.Range("I126").Value2 = 2525.5715
$ws.Range("J126").Value2 = 2518.3076
$ws.Range("K126").Value2 = 7576.7145
$ws.Range("L126").Value2 = 7554.9228
$ws.Range("M126").Value2 = -5106.7145
$ws.Range("N126").Value2 = -12494.9228
$ws.Range("H132").Value2 = 5684589
$ws.Range("I132").Value2 = 2751.6155
$ws.Range("J132").Value2 = 13891688
$ws.Range("K132").Value2 = 8254.8465
$ws.Range("L132").Value2 = 41675064
$ws.Range("M132").Value2 = -5724.8465
$ws.Range("N132").Value2 = -41680124
$ws.Range("H134").Value2 = 1242.1428
$ws.Range("I134").Value2 = 1242.1428
$ws.Range("J134").Value2 = 0
$ws.Range("K134").Value2 = 3726.4284
$ws.Range("L134").Value2 = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value2 = -1191.4284
$ws.Range("H136").Value2 = 2050
$ws.Range("I136").Value2 = 2300
$ws.Range("J136").Value2 = 1425
$ws.Range("K136").Value2 = 6900
$ws.Range("L136").Value2 = 4275
$ws.Range("M136").Value2 = -4350
$ws.Range("N136").Value2 = -9375
$ws.Range("H138").Value2 = 37933.332
$ws.Range("J138").Value2 = 37933.332
$ws.Range("L138").Value2 = 37933.332
$ws.Range("N138").Value2 = -48213.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value2 = 488.83334
$ws.Range("I68").Value2 = 435
$ws.Range("J68").Value2 = 596.5
$ws.Range("K68").Value2 = 1305
$ws.Range("L68").Value2 = 1789.5
$ws.Range("M68").Value2 = -494
$ws.Range("N68").Value2 = -3411.5
$ws.Range("H71").Value2 = 488.83334
$ws.Range("I71").Value2 = 435
$ws.Range("J71").Value2 = 596.5
$ws.Range("K71").Value2 = 3915
$ws.Range("L71").Value2 = 5368.5
$ws.Range("M71").Value2 = 141
$ws.Range("N71").Value2 = -13480.5
$ws.Range("H122").Value2 = 912.85
$ws.Range("I122").Value2 = 400.36365
$ws.Range("J122").Value2 = 1539.2222
$ws.Range("K122").Value2 = 3603.27285
$ws.Range("L122").Value2 = 13852.9998
$ws.Range("M122").Value2 = -1153.27285
$ws.Range("N122").Value2 = -18752.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 5028
$ws.Range("I80").Value2 = 3370.4
$ws.Range("J80").Value2 = 7100
$ws.Range("K80").Value2 = 3370.4
$ws.Range("L80").Value2 = 7100
$ws.Range("M80").Value2 = -2372.4
$ws.Range("N80").Value2 = -9096
$ws.Range("H83").Value2 = 5028
$ws.Range("I83").Value2 = 3370.4
$ws.Range("J83").Value2 = 7100
$ws.Range("K83").Value2 = 16852
$ws.Range("L83").Value2 = 35500
$ws.Range("M83").Value2 = -11860
$ws.Range("N83").Value2 = -45484
$ws.Range("H113").Value2 = 20834118
$ws.Range("I113").Value2 = 41667324
$ws.Range("K113").Value2 = 41667324
$ws.Range("M113").Value2 = -41665154
$ws.Range("H126").Value2 = 2759.375
$ws.Range("I126").Value2 = 3072.7273
$ws.Range("J126").Value2 = 2070
$ws.Range("K126").Value2 = 9218.1819
$ws.Range("L126").Value2 = 6210
$ws.Range("M126").Value2 = -6748.1819
$ws.Range("N126").Value2 = -11150
$ws.Range("H139").Value2 = 24000
$ws.Range("J139").Value2 = 24000
$ws.Range("L139").Value2 = 24000
$ws.Range("N139").Value2 = -34280
$ws.Range("H140").Value2 = 50000
$ws.Range("J140").Value2 = 50000
$ws.Range("L140").Value2 = 50000
$ws.Range("N140").Value2 = -60360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 799.6
$ws.Range("I16").Value2 = 349.5
$ws.Range("J16").Value2 = 1099.6666
$ws.Range("K16").Value2 = 349.5
$ws.Range("L16").Value2 = 1099.6666
$ws.Range("M16").Value2 = -179.5
$ws.Range("N16").Value2 = -1439.6666
$ws.Range("H133").Value2 = 22348.268
$ws.Range("J133").Value2 = 22348.268
$ws.Range("L133").Value2 = 22348.268
$ws.Range("N133").Value2 = -27408.268

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I17").Value2 = 2000
$ws.Range("J17").Value2 = 0
$ws.Range("K17").Value2 = 2000
$ws.Range("L17").ClearContents()
$ws.Range("M17").Value2 = -1828
$ws.Range("N17").Value2 = 0
$ws.Range("H64").Value2 = 19740
$ws.Range("J64").Value2 = 19740
$ws.Range("L64").Value2 = 19740
$ws.Range("N64").Value2 = -20236
$ws.Range("H67").Value2 = 19740
$ws.Range("J67").Value2 = 19740
$ws.Range("L67").Value2 = 19740
$ws.Range("N67").Value2 = -21456
$ws.Range("H136").Value2 = 3069.5557
$ws.Range("I136").Value2 = 3026.077
$ws.Range("J136").Value2 = 4200
$ws.Range("K136").Value2 = 9078.231
$ws.Range("L136").Value2 = 12600
$ws.Range("M136").Value2 = -6528.231
$ws.Range("N136").Value2 = -17700
